$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.683.98"
$ws.Range("E2").Value = "  +2.92%  "
$ws.Range("D3").Value = "'2.003.37"
$ws.Range("E3").Value = "  +6.75%  "
$ws.Range("D4").Value = "'1.009"
$ws.Range("E4").Value = "  +0.57%  "
$ws.Range("D5").Value = "'328.19"
$ws.Range("E5").Value = "  +0.96%  "
$ws.Range("D6").Value = "'1.008"
$ws.Range("E6").Value = "  +0.52%  "
$ws.Range("D7").Value = "'0.4700"
$ws.Range("E7").Value = "  +2.11%  "
$ws.Range("D8").Value = "'0.3958"
$ws.Range("E8").Value = "  +2.01%  "
$ws.Range("D9").Value = "'0.07983"
$ws.Range("E10").Value = "  +2.28%  "
$ws.Range("D11").Value = "'22.94"
$ws.Range("E11").Value = "  +5.42%  "
$ws.Range("D12").Value = "'2.065.68"
$ws.Range("E12").Value = "  +6.01%  "
$ws.Range("D13").Value = "'7.279"
$ws.Range("E13").Value = "  +3.88%  "
$ws.Range("D14").Value = "'5.898"
$ws.Range("E14").Value = "  +4.14%  "
$ws.Range("D15").Value = "'0.07175"
$ws.Range("E15").Value = "  +3.12%  "
$ws.Range("D16").Value = "'89.24"
$ws.Range("E16").Value = "  +1.09%  "
$ws.Range("D17").Value = "'1.010"
$ws.Range("E17").Value = "  +0.62%  "
$ws.Range("D18").Value = "'0.00001002"
$ws.Range("E18").Value = "  +0.54%  "
$ws.Range("D19").Value = "'17.40"
$ws.Range("E19").Value = "  +2.64%  "
$ws.Range("D20").Value = "'1.008"
$ws.Range("E20").Value = "  +0.48%  "
$ws.Range("D21").Value = "'29.773.33"
$ws.Range("E21").Value = "  +3.15%  "
$ws.Range("D22").Value = "'5.557"
$ws.Range("E22").Value = "  +5.57%  "
$ws.Range("D23").Value = "'11.31"
$ws.Range("E23").Value = "  +3.37%  "
$ws.Range("D24").Value = "'2.252.83"
$ws.Range("E24").Value = "  +7.49%  "
$ws.Range("D25").Value = "'2.131"
$ws.Range("E25").Value = "  +2.18%  "
$ws.Range("D26").Value = "'159.08"
$ws.Range("E26").Value = "  +2.22%  "
$ws.Range("D27").Value = "'19.78"
$ws.Range("E27").Value = "  +2.53%  "
$ws.Range("D28").Value = "'5.971"
$ws.Range("E28").Value = "  -0.20%  "
$ws.Range("B29").Value = "BitcoinCash"
$ws.Range("C29").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D29").Value = "'120.79"
$ws.Range("E29").Value = "  +2.83%  "
$ws.Range("B30").Value = "LidoDAOToken"
$ws.Range("C30").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D30").Value = "'1.974"
$ws.Range("E30").Value = "  +2.19%  "
$ws.Range("D31").Value = "'0.09478"
$ws.Range("E31").Value = "  +1.45%  "
$ws.Range("D32").Value = "'0.8999"
$ws.Range("E32").Value = "  -0.49%  "
$ws.Range("D33").Value = "'5.306"
$ws.Range("E33").Value = "  +0.78%  "
$ws.Range("D34").Value = "'1.349"
$ws.Range("E34").Value = "  +1.91%  "
$ws.Range("D35").Value = "'3.202"
$ws.Range("E35").Value = "  -1.95%  "
$ws.Range("D36").Value = "'0.05853"
$ws.Range("E36").Value = "  +1.50%  "
$ws.Range("B37").Value = "TrustWalletToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D37").Value = "'1.181"
$ws.Range("E37").Value = "  -0.59%  "
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").Value = "'0.02141"
$ws.Range("E38").Value = "  +3.36%  "
$ws.Range("B39").Value = "FraxShare"
$ws.Range("C39").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D39").Value = "'7.946"
$ws.Range("E39").Value = "  +3.15%  "
$ws.Range("B40").Value = "TheSandbox"
$ws.Range("C40").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D40").Value = "'0.5783"
$ws.Range("E40").Value = "  +2.22%  "
$ws.Range("B41").Value = "PEPE"
$ws.Range("C41").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D41").Value = "'0.000003160"
$ws.Range("E41").Value = "  +92.80%  "
$ws.Range("D42").Value = "'0.1827"
$ws.Range("E42").Value = "  +3.36%  "
$ws.Range("D43").Value = "'9.885"
$ws.Range("E43").Value = "  +2.06%  "
$ws.Range("D44").Value = "'12.17"
$ws.Range("E44").Value = "  +2.63%  "
$ws.Range("D45").Value = "'0.5405"
$ws.Range("E45").Value = "  +1.03%  "
$ws.Range("B46").Value = "MXToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D46").Value = "'2.669"
$ws.Range("E46").Value = "  +6.32%  "
$ws.Range("B47").Value = "RenderToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D47").Value = "'2.169"
$ws.Range("E47").Value = "  -4.25%  "
$ws.Range("B48").Value = "Cronos"
$ws.Range("C48").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D48").Value = "'0.06991"
$ws.Range("E48").Value = "  -0.66%  "
$ws.Range("B49").Value = "NEARProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D49").Value = "'1.877"
$ws.Range("E49").Value = "  +1.64%  "
$ws.Range("D50").Value = "'114.88"
$ws.Range("E50").Value = "  +1.67%  "
$ws.Range("E51").Value = "  +9.75%  "
